$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (User 5621813206 / Saad): update Last Activity, add Last Welcome
$ws.Range("G2").Value = "2025-11-12 11:31:26"
$ws.Range("H2").Value = "2025-11-12 11:00:43"

# Row 3 (User 8182750073 / Maher): change Language, update Last Updated, add Last Activity
$ws.Range("C3").Value = "ar"
$ws.Range("F3").Value = "2025-11-12 11:07:50"
$ws.Range("G3").Value = "2025-11-12 11:08:05"

$wb.Save()
